$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D for "Responsibility" (pushes NumberOfDays -> E, Tips -> F)
$ws.Columns.Item(4).Insert()

# Update header row (order chosen to reproduce shared-string table ordering)
$ws.Range("C1").Value = "Start Plus"
$ws.Range("B1").Value = 'Task "Title"'
$ws.Range("E1").Value = "Length (Days)"
$ws.Range("D1").Value = "Responsibility"

# Fill in the new Responsibility column with "EAO" for every task row
$ws.Range("D2").Value = "EAO"
$ws.Range("D3").Value = "EAO"
$ws.Range("D4").Value = "EAO"
$ws.Range("D5").Value = "EAO"
$ws.Range("D6").Value = "EAO"

# Column widths (best effort - engine quantizes ColumnWidth to whole pixels)
$ws.Columns.Item(2).ColumnWidth = 12.833333333333332
$ws.Columns.Item(4).ColumnWidth = 15.666666666666666
$ws.Columns.Item(5).ColumnWidth = 14.5
$ws.Columns.Item(6).ColumnWidth = 23.333333333333336

# Move the duplicate-values conditional formatting from the old Tips column (E) to the new one (F)
$fc = $ws.Range("E7:E51").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("F7:F51"))

# Restore the active cell/selection shown in the saved file
$ws.Range("E9").Select()
